$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New expense entry: row 3 ("Knowtefy DSC & DIN" paid in cash to Rabindra (CA))
$ws.Range("A3").Value = 2

# Date (20-Jul-2017) — paste the existing date cell's format so the new
# cell reuses the same numeric date style instead of minting a new one.
$ws.Range("B3").Value = 42936
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Range("C3").Value = "Knowtefy DSC & DIN"
$ws.Range("D3").Value = "Rabindra (CA)"
$ws.Range("E3").Value = 4000
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 4000
$ws.Range("H3").Value = "Cash Payment"

# Column B now holds dates, so let it size itself to the new content.
$ws.Columns.Item(2).EntireColumn.AutoFit()

# Leave the selection where the user last navigated to.
[void]$ws.Range("C15").Select()
